# Applies the two changes described by the commit's diff:
#
#   1. The table on slide 6 gets a new table style applied
#      (tableStyleId {8466050B-...} -> {B8477951-...}).
#
#   2. The colour palette of the deck's (single) theme -- the one actually
#      driving the slide master / all the slides -- is swapped from the
#      "Integral" palette to the stock "Office Theme" palette (the palette
#      that, in this deck, previously only lived -- unused -- on the Notes
#      Master's theme part).
#
# NB: PowerPoint's object model has no call that reassigns which on-disk
# theme part backs the slide master vs. the notes master, nor one that
# renames a theme/colour-scheme; the only programmatically reachable,
# persisted effect of "switch the active design from Integral to Office
# Theme" is rewriting the twelve RGB slots of the active ThemeColorScheme,
# which is what we do below.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 6.
# ---------------------------------------------------------------------
$targetStyleId = "{B8477951-69E1-4936-9A4C-AE54DEFF169C}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colours: Integral -> Office Theme on the active design.
# ---------------------------------------------------------------------
function Set-ThemeColor {
    param($scheme, [int]$index, [int]$r, [int]$g, [int]$b)
    # COM RGB longs are packed 0x00BBGGRR.
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $scheme.Item($index).RGB = $bgr
}

$design = $p.Designs.Item(1)
$colorScheme = $design.SlideMaster.Theme.ThemeColorScheme

Set-ThemeColor $colorScheme 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $colorScheme 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $colorScheme 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $colorScheme 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $colorScheme 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $colorScheme 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $colorScheme 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $colorScheme 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $colorScheme 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $colorScheme 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $colorScheme 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $colorScheme 12 0x95 0x4F 0x72   # folHlink

Write-Output "Applied table style $targetStyleId and Office Theme colours to the active design."
